$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bracingAssignment")
$ws.Range("B3:B66").Value = "V1"
$ws.Range("F60").Select()
